$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row-level updates: row -> @{ Col = NewValue; ... }
# Matches the "Updated cryptos list" GitHub Actions commit:
#   prices/volumes refreshed, and a handful of rows reordered
#   (TrustWalletToken/FraxShare, Aave/BabyDogeCoin, Mantle/Cronos).
$updates = [ordered]@{
    2 = @{ D='26.625.04'; E='  -7.35%  ' }
    3 = @{ D='1.696.86'; E='  -5.99%  ' }
    4 = @{ E='  +0.06%  ' }
    5 = @{ D='218.68'; E='  -5.82%  ' }
    6 = @{ D='0.5090'; E='  -14.09%  ' }
    7 = @{ D='1.003'; E='  -0.07%  ' }
    8 = @{ D='0.2642'; E='  -4.96%  ' }
    9 = @{ D='22.12'; E='  -5.26%  ' }
    10 = @{ D='0.06249'; E='  -8.55%  ' }
    11 = @{ D='0.07304'; E='  -2.62%  ' }
    12 = @{ D='1.690.57'; E='  -6.36%  ' }
    13 = @{ D='4.480'; E='  -6.20%  ' }
    14 = @{ D='0.5820'; E='  -6.80%  ' }
    15 = @{ D='1.926.15'; E='  -6.06%  ' }
    16 = @{ D='0.000008321'; E='  -10.53%  ' }
    17 = @{ D='65.34'; E='  -13.79%  ' }
    18 = @{ D='26.662.71'; E='  -7.01%  ' }
    19 = @{ D='5.041'; E='  -8.10%  ' }
    20 = @{ E='  -0.02%  ' }
    21 = @{ D='10.84'; E='  -5.57%  ' }
    22 = @{ D='186.31'; E='  -11.76%  ' }
    23 = @{ D='6.239'; E='  -8.96%  ' }
    24 = @{ E='  -0.04%  ' }
    25 = @{ D='144.95'; E='  -6.09%  ' }
    26 = @{ D='7.577'; E='  -4.00%  ' }
    27 = @{ D='0.1143'; E='  -10.10%  ' }
    28 = @{ D='15.66'; E='  -4.75%  ' }
    29 = @{ D='1.297'; E='  -8.79%  ' }
    30 = @{ D='0.05711'; E='  -7.86%  ' }
    31 = @{ D='1.329'; E='  -6.66%  ' }
    32 = @{ D='3.510'; E='  -6.76%  ' }
    33 = @{ D='3.494'; E='  -7.73%  ' }
    34 = @{ D='1.657'; E='  -4.16%  ' }
    35 = @{ D='1.017'; E='  -4.56%  ' }
    36 = @{ D='0.5968'; E='  -7.08%  ' }
    37 = @{ D='2.376'; E='  -4.70%  ' }
    38 = @{ D='2.683'; E='  -1.42%  ' }
    39 = @{ D='0.01598'; E='  -6.81%  ' }
    40 = @{ D='1.082.57'; E='  -5.46%  ' }
    41 = @{ B='FraxShare'; C='https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'; D='5.916'; E='  -9.92%  ' }
    42 = @{ B='TrustWalletToken'; C='https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'; D='0.8634'; E='  -2.01%  ' }
    43 = @{ D='1.002'; E='  -0.49%  ' }
    44 = @{ D='98.46'; E='  -1.86%  ' }
    45 = @{ D='1.854.26'; E='  -5.45%  ' }
    46 = @{ B='Aave'; C='https://coinranking.com/coin/ixgUfzmLR+aave-aave'; D='56.70'; E='  -6.34%  ' }
    47 = @{ B='BabyDogeCoin'; C='https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'; D='0.00000000107'; E='  -4.88%  ' }
    48 = @{ D='1.005'; E='  +0.03%  ' }
    49 = @{ D='8.116'; E='  -3.30%  ' }
    50 = @{ B='Cronos'; C='https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'; D='0.05238'; E='  -4.21%  ' }
    51 = @{ B='Mantle'; C='https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'; D='0.4315'; E='  -3.80%  ' }
}

foreach ($row in $updates.Keys) {
    $rowData = $updates[$row]
    foreach ($col in $rowData.Keys) {
        $addr = "$col$row"
        $cell = $ws.Range($addr)
        if ($col -eq "D" -or $col -eq "E") {
            # Force text format so numeric-looking strings (prices like
            # "1.690.57" or volumes like "  -6.36%  ") are not coerced
            # into actual Double values by Excel.
            $cell.NumberFormat = "@"
        }
        $cell.Value = $rowData[$col]
    }
}
